$d = $word.ActiveDocument

function Set-ParagraphXml($FindText, $ParagraphInnerXml) {
    $rng = $d.Content
    $rng.Find.Execute($FindText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $para = $rng.Paragraphs(1)
    $fullRange = $para.Range

    $pkg = '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $ParagraphInnerXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $fullRange.InsertXML($pkg)
}

# 1) "Honor's Prog" + bookmark + "ram Member" -> single run "Honor's Program Member",
#    dropping the _GoBack bookmark from here.
Set-ParagraphXml "Honor" (
    '<w:p w:rsidR="00B07427" w:rsidRDefault="005932E6" w:rsidP="000B005F">' +
    '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr>' +
    '<w:spacing w:line="240" w:lineRule="auto"/><w:contextualSpacing/></w:pPr>' +
    '<w:r><w:t>Honor’s Program Member</w:t></w:r>' +
    '</w:p>'
)

# 2) Add ", FXML" run after ", C, HTML5, CSS3"
Set-ParagraphXml "HTML5, CSS3" (
    '<w:p w:rsidR="00B07427" w:rsidRDefault="005932E6" w:rsidP="000B005F">' +
    '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr>' +
    '<w:spacing w:line="240" w:lineRule="auto"/><w:contextualSpacing/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Java, C++, ARM assembly, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Javascript</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">, JQuery, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Labview</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>, C, HTML5, CSS3</w:t></w:r>' +
    '<w:r><w:t>, FXML</w:t></w:r>' +
    '</w:p>'
)

# 3) Add ", Netbeans, Eclipse, Unity development" after ", Autodesk Suite"
Set-ParagraphXml "Autodesk Suite" (
    '<w:p w:rsidR="00B07427" w:rsidRDefault="005932E6" w:rsidP="000B005F">' +
    '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr>' +
    '<w:spacing w:line="240" w:lineRule="auto"/><w:contextualSpacing/></w:pPr>' +
    '<w:r><w:t>Windows, Linux, Word/Excel/</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Powerpoint</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>, Android dev</w:t></w:r>' +
    '<w:r w:rsidR="000B005F"><w:t>elopment</w:t></w:r>' +
    '<w:r><w:t>, Autodesk Suite</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Netbeans</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>, Eclipse, Unity development</w:t></w:r>' +
    '</w:p>'
)

# 4) Move the _GoBack bookmark to the end of the "Foreign Languages: Fluent French" paragraph
Set-ParagraphXml "Foreign Languages" (
    '<w:p w:rsidR="00B07427" w:rsidRDefault="005932E6" w:rsidP="000B005F">' +
    '<w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>' +
    '<w:r w:rsidRPr="000B005F"><w:rPr><w:b/></w:rPr><w:t>Foreign Languages</w:t></w:r>' +
    '<w:r><w:t>: Fluent French</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'
)
